$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for "Feria Lagunitas de Puerto
# Montt - Cebollín". It belongs chronologically right before the row that
# used to be row 102 (2021-08-06), so insert a fresh row at 102 and push
# everything from the old row 102 down to row 103 (...through old row 159,
# which becomes row 160).
$ws.Rows.Item(102).Insert()

# Columns that are constant for every record in this per-product sheet.
$ws.Cells.Item(102, 1).Value = 4
$ws.Cells.Item(102, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(102, 3).Value = "Los Lagos"
$ws.Cells.Item(102, 5).Value = 10
$ws.Cells.Item(102, 6).Value = 100112037
$ws.Cells.Item(102, 7).Value = "Cebollín"
$ws.Cells.Item(102, 8).Value = "Sin especificar"
$ws.Cells.Item(102, 9).Value = "Primera"
$ws.Cells.Item(102, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(102, 15).Value = "Región Metropolitana"
$ws.Cells.Item(102, 17).Value = 36
$ws.Cells.Item(102, 18).Value = "Hortaliza"

# New row-specific observation values.
$ws.Cells.Item(102, 4).Value = 44460
$ws.Cells.Item(102, 10).Value = 90
$ws.Cells.Item(102, 11).Value = 6000
$ws.Cells.Item(102, 12).Value = 6000
$ws.Cells.Item(102, 13).Value = 6000
$ws.Cells.Item(102, 16).Value = 167
